$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("B9").Value = 6834733
$ws.Range("C9").Value = "Croatia 1NL"
$ws.Range("D9").Value = 45150.52083333334
$ws.Range("E9").Value = "HNK Cibalia"
$ws.Range("F9").Value = "NK Croatia Zmijavci"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = "H"
$ws.Range("L9").Value = 1.65
$ws.Range("M9").Value = 3.5
$ws.Range("N9").Value = 4.5
$ws.Range("O9").Value = 1.909
$ws.Range("P9").Value = 3.3
$ws.Range("Q9").Value = 3.3
$ws.Range("R9").Value = -0.5
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 1.8
$ws.Range("U9").Value = 2.25
$ws.Range("V9").Value = 1.95
$ws.Range("W9").Value = 1.85
$ws.Range("X9").Value = 0.909
$ws.Range("Y9").Value = -1
$ws.Range("Z9").Value = -1
$ws.Range("AA9").Value = 1
$ws.Range("AB9").Value = -1
$ws.Range("AC9").Value = -1
$ws.Range("AD9").Value = 0.8500000000000001

# Row 10
$ws.Range("B10").Value = 6834729
$ws.Range("C10").Value = "Croatia 1NL"
$ws.Range("D10").Value = 45150.52083333334
$ws.Range("E10").Value = "NK Solin"
$ws.Range("F10").Value = "Bijelo Brdo"
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = "D"
$ws.Range("L10").Value = 2.1
$ws.Range("M10").Value = 3.2
$ws.Range("N10").Value = 3.1
$ws.Range("O10").Value = 2.05
$ws.Range("P10").Value = 3.25
$ws.Range("Q10").Value = 3.25
$ws.Range("R10").Value = -0.25
$ws.Range("S10").Value = 1.8
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 2.5
$ws.Range("V10").Value = 2
$ws.Range("W10").Value = 1.8
$ws.Range("X10").Value = -1
$ws.Range("Y10").Value = 2.25
$ws.Range("Z10").Value = -1
$ws.Range("AA10").Value = -0.5
$ws.Range("AB10").Value = 0.5
$ws.Range("AC10").Value = -1
$ws.Range("AD10").Value = 0.8

# Row 11
$ws.Range("B11").Value = 6834732
$ws.Range("C11").Value = "Croatia 1NL"
$ws.Range("D11").Value = 45150.52083333334
$ws.Range("E11").Value = "Orijent"
$ws.Range("F11").Value = "NK Dubrava Zagreb"
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = "D"
$ws.Range("L11").Value = 2.1
$ws.Range("M11").Value = 3.2
$ws.Range("N11").Value = 3.1
$ws.Range("O11").Value = 2.3
$ws.Range("P11").Value = 3.1
$ws.Range("Q11").Value = 2.8
$ws.Range("R11").Value = -0.25
$ws.Range("S11").Value = 2.025
$ws.Range("T11").Value = 1.775
$ws.Range("U11").Value = 2.25
$ws.Range("V11").Value = 1.825
$ws.Range("W11").Value = 1.975
$ws.Range("X11").Value = -1
$ws.Range("Y11").Value = 2.1
$ws.Range("Z11").Value = -1
$ws.Range("AA11").Value = -0.5
$ws.Range("AB11").Value = 0.3875
$ws.Range("AC11").Value = -0.5
$ws.Range("AD11").Value = 0.4875

# Row 21
$ws.Range("B21").Value = 6834743
$ws.Range("C21").Value = "Croatia 1NL"
$ws.Range("D21").Value = 45164.5
$ws.Range("E21").Value = "Orijent"
$ws.Range("F21").Value = "Vukovar 91"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = "A"
$ws.Range("L21").Value = 3
$ws.Range("M21").Value = 3.3
$ws.Range("N21").Value = 2.1
$ws.Range("O21").Value = 1.95
$ws.Range("P21").Value = 3.4
$ws.Range("Q21").Value = 3.3
$ws.Range("R21").Value = -0.5
$ws.Range("S21").Value = 2
$ws.Range("T21").Value = 1.8
$ws.Range("U21").Value = 2.5
$ws.Range("V21").Value = 2
$ws.Range("W21").Value = 1.8
$ws.Range("X21").Value = -1
$ws.Range("Y21").Value = -1
$ws.Range("Z21").Value = 2.3
$ws.Range("AA21").Value = -1
$ws.Range("AB21").Value = 0.8
$ws.Range("AC21").Value = -1
$ws.Range("AD21").Value = 0.8

# Row 23
$ws.Range("B23").Value = 6834744
$ws.Range("C23").Value = "Croatia 1NL"
$ws.Range("D23").Value = 45164.5
$ws.Range("E23").Value = "HNK Cibalia"
$ws.Range("F23").Value = "NK Jarun"
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = "H"
$ws.Range("L23").Value = 1.75
$ws.Range("M23").Value = 3.4
$ws.Range("N23").Value = 4
$ws.Range("O23").Value = 1.666
$ws.Range("P23").Value = 3.6
$ws.Range("Q23").Value = 4.2
$ws.Range("R23").Value = -0.75
$ws.Range("S23").Value = 1.925
$ws.Range("T23").Value = 1.875
$ws.Range("U23").Value = 2.5
$ws.Range("V23").Value = 1.825
$ws.Range("W23").Value = 1.975
$ws.Range("X23").Value = 0.6659999999999999
$ws.Range("Y23").Value = -1
$ws.Range("Z23").Value = -1
$ws.Range("AA23").Value = 0.925
$ws.Range("AB23").Value = -1
$ws.Range("AC23").Value = 0.825
$ws.Range("AD23").Value = -1

# Row 142
$ws.Range("B142").Value = 7977238
$ws.Range("C142").Value = "Croatia 1NL"
$ws.Range("D142").Value = 45374.47916666666
$ws.Range("E142").Value = "NK Jarun"
$ws.Range("F142").Value = "Bijelo Brdo"
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 3
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = 2
$ws.Range("K142").Value = "A"
$ws.Range("L142").Value = 1.909
$ws.Range("M142").Value = 3.3
$ws.Range("N142").Value = 3.5
$ws.Range("O142").Value = 1.909
$ws.Range("P142").Value = 3.4
$ws.Range("Q142").Value = 3.5
$ws.Range("R142").Value = -0.5
$ws.Range("S142").Value = 1.95
$ws.Range("T142").Value = 1.85
$ws.Range("U142").Value = 2.25
$ws.Range("V142").Value = 1.75
$ws.Range("W142").Value = 2.05
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = 2.5
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.8500000000000001
$ws.Range("AC142").Value = 0.75
$ws.Range("AD142").Value = -1

# Row 143
$ws.Range("B143").Value = 7977239
$ws.Range("C143").Value = "Croatia 1NL"
$ws.Range("D143").Value = 45374.47916666666
$ws.Range("E143").Value = "NK Dugopolje"
$ws.Range("F143").Value = "NK Solin"
$ws.Range("G143").Value = 2
$ws.Range("H143").Value = 2
$ws.Range("I143").Value = 0
$ws.Range("J143").Value = 1
$ws.Range("K143").Value = "D"
$ws.Range("L143").Value = 2
$ws.Range("M143").Value = 3.2
$ws.Range("N143").Value = 3.3
$ws.Range("O143").Value = 2.3
$ws.Range("P143").Value = 3.2
$ws.Range("Q143").Value = 2.7
$ws.Range("R143").Value = -0.25
$ws.Range("S143").Value = 1.975
$ws.Range("T143").Value = 1.725
$ws.Range("U143").Value = 2.25
$ws.Range("V143").Value = 1.9
$ws.Range("W143").Value = 1.9
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 2.2
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = -0.5
$ws.Range("AB143").Value = 0.3625
$ws.Range("AC143").Value = 0.8999999999999999
$ws.Range("AD143").Value = -1

# Row 203
$ws.Range("B203").Value = 7977296
$ws.Range("C203").Value = "Croatia 1NL"
$ws.Range("D203").Value = 45444.52083333334
$ws.Range("E203").Value = "HNK Cibalia"
$ws.Range("F203").Value = "Vukovar 91"
$ws.Range("G203").Value = 1
$ws.Range("H203").Value = 5
$ws.Range("K203").Value = "A"
$ws.Range("L203").Value = 2.75
$ws.Range("M203").Value = 3.4
$ws.Range("N203").Value = 2.2
$ws.Range("O203").Value = 3
$ws.Range("P203").Value = 3.6
$ws.Range("Q203").Value = 2.05
$ws.Range("R203").Value = 0.25
$ws.Range("S203").Value = 1.975
$ws.Range("T203").Value = 1.825
$ws.Range("U203").Value = 2.75
$ws.Range("V203").Value = 1.85
$ws.Range("W203").Value = 1.95
$ws.Range("X203").Value = -1
$ws.Range("Y203").Value = -1
$ws.Range("Z203").Value = 1.05
$ws.Range("AA203").Value = -1
$ws.Range("AB203").Value = 0.825
$ws.Range("AC203").Value = 0.8500000000000001
$ws.Range("AD203").Value = -1

# Row 204
$ws.Range("B204").Value = 7977297
$ws.Range("C204").Value = "Croatia 1NL"
$ws.Range("D204").Value = 45444.52083333334
$ws.Range("E204").Value = "NK Croatia Zmijavci"
$ws.Range("F204").Value = "NK Sesvete"
$ws.Range("G204").Value = 3
$ws.Range("H204").Value = 1
$ws.Range("K204").Value = "H"
$ws.Range("L204").Value = 1.5
$ws.Range("M204").Value = 3.4
$ws.Range("N204").Value = 6.5
$ws.Range("O204").Value = 1.4
$ws.Range("P204").Value = 3.8
$ws.Range("Q204").Value = 6.5
$ws.Range("R204").Value = -1.25
$ws.Range("S204").Value = 1.925
$ws.Range("T204").Value = 1.875
$ws.Range("U204").Value = 3
$ws.Range("V204").Value = 2
$ws.Range("W204").Value = 1.8
$ws.Range("X204").Value = 0.3999999999999999
$ws.Range("Y204").Value = -1
$ws.Range("Z204").Value = -1
$ws.Range("AA204").Value = 0.925
$ws.Range("AB204").Value = -1
$ws.Range("AC204").Value = 1
$ws.Range("AD204").Value = -1
